$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 116
$ws.Range("H116").Value = 6503
$ws.Range("I116").Value = 7942.625
$ws.Range("J116").Value = 3623.75
$ws.Range("K116").Value = 7942.625
$ws.Range("L116").Value = 3623.75
$ws.Range("M116").Value = -4500.625
$ws.Range("N116").Value = -10507.75

# ALC row 129
$ws.Range("H129").Value = 1017.44684
$ws.Range("J129").Value = 1168.5143
$ws.Range("L129").Value = 3505.5429
$ws.Range("N129").Value = -13505.5429

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 871.7308
$ws.Range("I2").Value = 840.2083
$ws.Range("J2").Value = 1250
$ws.Range("K2").Value = 840.2083
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = -727.2083
$ws.Range("N2").Value = -1476

# ARM row 45
$ws.Range("H45").Value = 10850.637
$ws.Range("I45").Value = 10935.7
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 10935.7
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -10558.7
$ws.Range("N45").Value = -10754

# ARM row 61
$ws.Range("H61").Value = 2664.8096
$ws.Range("I61").Value = 2047.8334
$ws.Range("J61").Value = 6366.6665
$ws.Range("K61").Value = 2047.8334
$ws.Range("L61").Value = 6366.6665
$ws.Range("M61").Value = -1835.8334
$ws.Range("N61").Value = -6790.6665

# ARM row 116
$ws.Range("H116").Value = 871.7308
$ws.Range("I116").Value = 840.2083
$ws.Range("J116").Value = 1250
$ws.Range("K116").Value = 840.2083
$ws.Range("L116").Value = 1250
$ws.Range("M116").Value = 1453.7917
$ws.Range("N116").Value = -5838

# ARM row 136
$ws.Range("H136").Value = 2664.8096
$ws.Range("I136").Value = 2047.8334
$ws.Range("J136").Value = 6366.6665
$ws.Range("K136").Value = 6143.5002
$ws.Range("L136").Value = 19099.9995
$ws.Range("M136").Value = -3593.5002
$ws.Range("N136").Value = -24199.9995

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 871.7308
$ws.Range("I3").Value = 840.2083
$ws.Range("J3").Value = 1250
$ws.Range("K3").Value = 840.2083
$ws.Range("L3").Value = 1250
$ws.Range("M3").Value = -726.2083
$ws.Range("N3").Value = -1478

# BSM row 80
$ws.Range("H80").Value = 392.68182
$ws.Range("J80").Value = 438.57895
$ws.Range("L80").Value = 438.57895
$ws.Range("N80").Value = -2434.57895

# BSM row 83
$ws.Range("H83").Value = 392.68182
$ws.Range("J83").Value = 438.57895
$ws.Range("L83").Value = 2192.89475
$ws.Range("N83").Value = -12176.89475

# BSM row 134
$ws.Range("H134").Value = 15153744
$ws.Range("I134").Value = 20835294
$ws.Range("J134").Value = 2944.4443
$ws.Range("K134").Value = 62505882
$ws.Range("L134").Value = 8833.332900000001
$ws.Range("M134").Value = -62503347
$ws.Range("N134").Value = -13903.3329

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 1950.3334
$ws.Range("I16").Value = 1681
$ws.Range("J16").Value = 2287
$ws.Range("K16").Value = 1681
$ws.Range("L16").Value = 2287
$ws.Range("M16").Value = -1394
$ws.Range("N16").Value = -2861

# CRP row 31
$ws.Range("H31").Value = 4246.8984
$ws.Range("I31").Value = 1636.561
$ws.Range("J31").Value = 10192.667
$ws.Range("K31").Value = 1636.561
$ws.Range("L31").Value = 10192.667
$ws.Range("M31").Value = -1341.561
$ws.Range("N31").Value = -10782.667

# CRP row 34
$ws.Range("H34").Value = 4246.8984
$ws.Range("I34").Value = 1636.561
$ws.Range("J34").Value = 10192.667
$ws.Range("K34").Value = 1636.561
$ws.Range("L34").Value = 10192.667
$ws.Range("M34").Value = -1434.561
$ws.Range("N34").Value = -10596.667

# CRP row 58
$ws.Range("H58").Value = 1898.08
$ws.Range("I58").Value = 1363.3846
$ws.Range("J58").Value = 2477.3333
$ws.Range("K58").Value = 1363.3846
$ws.Range("L58").Value = 2477.3333
$ws.Range("M58").Value = -1160.3846
$ws.Range("N58").Value = -2883.3333

# CRP row 105
$ws.Range("H105").Value = 1703.2084
$ws.Range("I105").Value = 1914.25
$ws.Range("K105").Value = 1914.25
$ws.Range("M105").Value = -167.25

# CRP row 107
$ws.Range("H107").Value = 12346485
$ws.Range("I107").Value = 15152050
$ws.Range("J107").Value = 1999.8
$ws.Range("K107").Value = 15152050
$ws.Range("L107").Value = 1999.8
$ws.Range("M107").Value = -15150130
$ws.Range("N107").Value = -5839.8

# CRP row 113
$ws.Range("H113").Value = 1950.3334
$ws.Range("I113").Value = 1681
$ws.Range("J113").Value = 2287
$ws.Range("K113").Value = 1681
$ws.Range("L113").Value = 2287
$ws.Range("M113").Value = 489
$ws.Range("N113").Value = -6627

# CRP row 122
$ws.Range("H122").Value = 1391052.4
$ws.Range("I122").Value = 2527425.8
$ws.Range("J122").Value = 2151.6667
$ws.Range("K122").Value = 7582277.399999999
$ws.Range("L122").Value = 6455.000100000001
$ws.Range("M122").Value = -7579827.399999999
$ws.Range("N122").Value = -11355.0001

# CRP row 132
$ws.Range("H132").Value = 2220.32
$ws.Range("I132").Value = 1706.6666
$ws.Range("J132").Value = 3541.1428
$ws.Range("K132").Value = 5119.9998
$ws.Range("L132").Value = 10623.4284
$ws.Range("M132").Value = -2589.9998
$ws.Range("N132").Value = -15683.4284

# CRP row 133
$ws.Range("H133").Value = 30326
$ws.Range("J133").Value = 30326
$ws.Range("L133").Value = 30326
$ws.Range("N133").Value = -35386

# CRP row 134
$ws.Range("H134").Value = 2577.0715
$ws.Range("I134").Value = 2572.6216
$ws.Range("J134").Value = 2610
$ws.Range("K134").Value = 7717.864799999999
$ws.Range("L134").Value = 7830
$ws.Range("M134").Value = -5182.864799999999
$ws.Range("N134").Value = -12900

# CRP row 136
$ws.Range("H136").Value = 1898.08
$ws.Range("I136").Value = 1363.3846
$ws.Range("J136").Value = 2477.3333
$ws.Range("K136").Value = 4090.1538
$ws.Range("L136").Value = 7431.999899999999
$ws.Range("M136").Value = -1540.1538
$ws.Range("N136").Value = -12531.9999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 114
$ws.Range("H114").Value = 5590.7827
$ws.Range("I114").Value = 433.58334
$ws.Range("J114").Value = 11216.818
$ws.Range("K114").Value = 1300.75002
$ws.Range("L114").Value = 33650.454
$ws.Range("M114").Value = 1953.24998
$ws.Range("N114").Value = -40158.454

# CUL row 131
$ws.Range("H131").Value = 910.3099999999999
$ws.Range("I131").Value = 580.1111
$ws.Range("J131").Value = 942.96704
$ws.Range("K131").Value = 1740.3333
$ws.Range("L131").Value = 2828.90112
$ws.Range("M131").Value = 3299.6667
$ws.Range("N131").Value = -12908.90112

$ws = $wb.Worksheets.Item("GSM")
# GSM row 93
$ws.Range("H93").Value = 9250.833000000001
$ws.Range("J93").Value = 9250.833000000001
$ws.Range("L93").Value = 9250.833000000001
$ws.Range("N93").Value = -12994.833

# GSM row 126
$ws.Range("H126").Value = 6877
$ws.Range("I126").Value = 10100.167
$ws.Range("J126").Value = 2579.4443
$ws.Range("K126").Value = 30300.501
$ws.Range("L126").Value = 7738.3329
$ws.Range("M126").Value = -27830.501
$ws.Range("N126").Value = -12678.3329

# GSM row 132
$ws.Range("H132").Value = 5581
$ws.Range("I132").Value = 8527.200000000001
$ws.Range("J132").Value = 3944.2222
$ws.Range("K132").Value = 25581.6
$ws.Range("L132").Value = 11832.6666
$ws.Range("M132").Value = -23051.6
$ws.Range("N132").Value = -16892.6666

# GSM row 133
$ws.Range("H133").Value = 39890
$ws.Range("J133").Value = 39890
$ws.Range("L133").Value = 39890
$ws.Range("N133").Value = -50010

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 1623
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1623
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1623
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2213

# LTW row 27
$ws.Range("H27").Value = 1623
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1623
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1623
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1837

# LTW row 46
$ws.Range("H46").Value = 616.5
$ws.Range("I46").Value = 457.14285
$ws.Range("K46").Value = 457.14285
$ws.Range("M46").Value = -269.14285

$ws = $wb.Worksheets.Item("WVR")
# WVR row 126
$ws.Range("H126").Value = 1020.9231
$ws.Range("I126").Value = 753.55554
$ws.Range("J126").Value = 1622.5
$ws.Range("K126").Value = 2260.66662
$ws.Range("L126").Value = 4867.5
$ws.Range("M126").Value = 209.33338
$ws.Range("N126").Value = -9807.5

# WVR row 136
$ws.Range("H136").Value = 5175.115
$ws.Range("I136").Value = 8396.083000000001
$ws.Range("J136").Value = 2414.2856
$ws.Range("K136").Value = 25188.249
$ws.Range("L136").Value = 7242.8568
$ws.Range("M136").Value = -22638.249
$ws.Range("N136").Value = -12342.8568
